{"js": "const body = context.document.body;\nconst pairs = [\n  [\"237\u00f78=29, 5\", \"221\u00f75=44, 1\"],\n  [\"233\u00f74=58, 1\", \"176\u00f72=88, 0\"],\n  [\"978\u00f79=108, 6\", \"490\u00f79=54, 4\"],\n  [\"455\u00f75=91, 0\", \"665\u00f79=73, 8\"],\n  [\"181\u00f73=60, 1\", \"768\u00f74=192, 0\"],\n  [\"128\u00f76=21, 2\", \"728\u00f76=121, 2\"],\n  [\"342\u00f76=57, 0\", \"547\u00f76=91, 1\"],\n  [\"355\u00f72=177, 1\", \"608\u00f77=86, 6\"],\n  [\"279\u00f73=93, 0\", \"937\u00f77=133, 6\"],\n  [\"934\u00f75=186, 4\", \"692\u00f79=76, 8\"],\n  [\"422\u00f76=70, 2\", \"783\u00f77=111, 6\"],\n  [\"914\u00f72=457, 0\", \"754\u00f74=188, 2\"],\n  [\"203\u00f76=33, 5\", \"230\u00f78=28, 6\"],\n  [\"695\u00f79=77, 2\", \"407\u00f79=45, 2\"],\n  [\"304\u00f72=152, 0\", \"186\u00f72=93, 0\"],\n  [\"498\u00f76=83, 0\", \"863\u00f75=172, 3\"],\n  [\"713\u00f73=237, 2\", \"939\u00f73=313, 0\"],\n  [\"447\u00f79=49, 6\", \"892\u00f78=111, 4\"],\n  [\"718\u00f75=143, 3\", \"982\u00f75=196, 2\"],\n  [\"228\u00f73=76, 0\", \"228\u00f74=57, 0\"],\n  [\"360\u00f78=45, 0\", \"209\u00f72=104, 1\"],\n  [\"239\u00f78=29, 7\", \"259\u00f72=129, 1\"],\n  [\"996\u00f79=110, 6\", \"456\u00f75=91, 1\"],\n  [\"912\u00f76=152, 0\", \"532\u00f74=133, 0\"],\n  [\"432\u00f77=61, 5\", \"256\u00f77=36, 4\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"237\u00f78=29, 5\", \"221\u00f75=44, 1\"),\n    @(\"233\u00f74=58, 1\", \"176\u00f72=88, 0\"),\n    @(\"978\u00f79=108, 6\", \"490\u00f79=54, 4\"),\n    @(\"455\u00f75=91, 0\", \"665\u00f79=73, 8\"),\n    @(\"181\u00f73=60, 1\", \"768\u00f74=192, 0\"),\n    @(\"128\u00f76=21, 2\", \"728\u00f76=121, 2\"),\n    @(\"342\u00f76=57, 0\", \"547\u00f76=91, 1\"),\n    @(\"355\u00f72=177, 1\", \"608\u00f77=86, 6\"),\n    @(\"279\u00f73=93, 0\", \"937\u00f77=133, 6\"),\n    @(\"934\u00f75=186, 4\", \"692\u00f79=76, 8\"),\n    @(\"422\u00f76=70, 2\", \"783\u00f77=111, 6\"),\n    @(\"914\u00f72=457, 0\", \"754\u00f74=188, 2\"),\n    @(\"203\u00f76=33, 5\", \"230\u00f78=28, 6\"),\n    @(\"695\u00f79=77, 2\", \"407\u00f79=45, 2\"),\n    @(\"304\u00f72=152, 0\", \"186\u00f72=93, 0\"),\n    @(\"498\u00f76=83, 0\", \"863\u00f75=172, 3\"),\n    @(\"713\u00f73=237, 2\", \"939\u00f73=313, 0\"),\n    @(\"447\u00f79=49, 6\", \"892\u00f78=111, 4\"),\n    @(\"718\u00f75=143, 3\", \"982\u00f75=196, 2\"),\n    @(\"228\u00f73=76, 0\", \"228\u00f74=57, 0\"),\n    @(\"360\u00f78=45, 0\", \"209\u00f72=104, 1\"),\n    @(\"239\u00f78=29, 7\", \"259\u00f72=129, 1\"),\n    @(\"996\u00f79=110, 6\", \"456\u00f75=91, 1\"),\n    @(\"912\u00f76=152, 0\", \"532\u00f74=133, 0\"),\n    @(\"432\u00f77=61, 5\", \"256\u00f77=36, 4\"),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $ok = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $ok) {\n        throw \"Not found: $oldText\"\n    }\n}"}
